$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Make the shared-formula groups in the "gas" sheet explicit for
#    columns H, I, J, K (rows 13:19) -- matches the saved workbook's
#    shared-formula restructuring.
# ------------------------------------------------------------------
$gas = $wb.Worksheets.Item("gas")
$gas.Range("H13:H19").Formula = '=$B$2*($E13*100-$B$3)'
$gas.Range("I13:I19").Formula = '=$B$5*($E13*100-$B$6)'
$gas.Range("J13:J19").Formula = '=$I13-$B$7*($E13*100-$B$6)*($B13-$D$10)'
$gas.Range("K13:K19").Formula = '=$I13-$B$8*($E13*100-$B$6)*($B13-$D$10)'

# ------------------------------------------------------------------
# 2. Add the new "densification" sheet after "gas" (becomes the new
#    active / selected tab, same as in the target workbook).
# ------------------------------------------------------------------
$ws = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws.Name = "densification"

# -- inputs --------------------------------------------------------
$ws.Range("F4").Value = "total dens"
$ws.Range("G4").Value = -0.034
$ws.Range("F5").Value = "Bmax"
$ws.Range("G5").Value = 0.006

# -- header row (row 7) ---------------------------------------------
$ws.Range("C7").Value = "#"
$ws.Range("D7").Value = "Step"
$ws.Range("E7").Value = "burnup"
# "dens" must be interned (as a shared string) before "BUCK vol" so the
# shared-string table ends up in the same order as the target workbook.
$ws.Range("G7").Value = "dens"
$ws.Range("F7").Value = "BUCK vol"
$ws.Range("H7").Value = "EXCEL vol"
$ws.Range("I7").Value = "% Diff"

# -- data rows 8-13 --------------------------------------------------
$ws.Range("D8").Value = 10
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 1

$ws.Range("C9").Value = "3="
$ws.Range("D9").Value = 20
$ws.Range("E9").Value = 0.00365437
$ws.Range("F9").Value = 0.9844954

$ws.Range("D10").Value = 30
$ws.Range("E10").Value = 0.00730874
$ws.Range("F10").Value = 0.976062

$ws.Range("D11").Value = 40
$ws.Range("E11").Value = 0.01096311
$ws.Range("F11").Value = 0.9714751

$ws.Range("D12").Value = 50
$ws.Range("E12").Value = 0.01461748
$ws.Range("F12").Value = 0.9689804

$ws.Range("D13").Value = 60
$ws.Range("E13").Value = 0.01827185
$ws.Range("F13").Value = 0.9676236

# -- formulas in G:I, rows 8-13 (entered per-cell -- the source
#    workbook keeps these as independent, non-shared formulas) --------
$ws.Range("G8").Formula = '=$G$4*(1-EXP(-E8/$G$5))'
$ws.Range("H8").Formula = '=(1+G8)'
$ws.Range("I8").Formula = '=ABS(F8-H8)/H8*100'

$ws.Range("G9").Formula = '=$G$4*(1-EXP(-E9/$G$5))'
$ws.Range("H9").Formula = '=(1+G9)'
$ws.Range("I9").Formula = '=ABS(F9-H9)/H9*100'

$ws.Range("G10").Formula = '=$G$4*(1-EXP(-E10/$G$5))'
$ws.Range("H10").Formula = '=(1+G10)'
$ws.Range("I10").Formula = '=ABS(F10-H10)/H10*100'

$ws.Range("G11").Formula = '=$G$4*(1-EXP(-E11/$G$5))'
$ws.Range("H11").Formula = '=(1+G11)'
$ws.Range("I11").Formula = '=ABS(F11-H11)/H11*100'

$ws.Range("G12").Formula = '=$G$4*(1-EXP(-E12/$G$5))'
$ws.Range("H12").Formula = '=(1+G12)'
$ws.Range("I12").Formula = '=ABS(F12-H12)/H12*100'

$ws.Range("G13").Formula = '=$G$4*(1-EXP(-E13/$G$5))'
$ws.Range("H13").Formula = '=(1+G13)'
$ws.Range("I13").Formula = '=ABS(F13-H13)/H13*100'

# -- number format for rows 8-13, columns E:I (matches numFmt 166,
#    "0.0000E+00", reused as style index 7 in the source workbook) ----
$ws.Range("E8:I13").NumberFormat = "0.0000E+00"

# -- leftover formatted (empty) cells, rows 19-26 ---------------------
$ws.Range("G19:J24").NumberFormat = "0.00E+00"
$ws.Range("G25:I26").NumberFormat = "0.00E+00"

# -- selection matches the saved workbook ------------------------------
[void]$ws.Range("C7").Select()
